$d = $word.ActiveDocument

$d.Content.Find.Execute("2024-11-23 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-11-24 Sunday", 2) | Out-Null
$d.Content.Find.Execute("219×7=1533", $true, $false, $false, $false, $false, $true, 1, $false, "261×6=1566", 2) | Out-Null
$d.Content.Find.Execute("541×3=1623", $true, $false, $false, $false, $false, $true, 1, $false, "101×6=606", 2) | Out-Null
$d.Content.Find.Execute("166×4=664", $true, $false, $false, $false, $false, $true, 1, $false, "677×4=2708", 2) | Out-Null
$d.Content.Find.Execute("250×9=2250", $true, $false, $false, $false, $false, $true, 1, $false, "661×6=3966", 2) | Out-Null
$d.Content.Find.Execute("805×9=7245", $true, $false, $false, $false, $false, $true, 1, $false, "276×7=1932", 2) | Out-Null
$d.Content.Find.Execute("986×8=7888", $true, $false, $false, $false, $false, $true, 1, $false, "220×6=1320", 2) | Out-Null
$d.Content.Find.Execute("345×9=3105", $true, $false, $false, $false, $false, $true, 1, $false, "658×2=1316", 2) | Out-Null
$d.Content.Find.Execute("222×2=444", $true, $false, $false, $false, $false, $true, 1, $false, "607×5=3035", 2) | Out-Null
$d.Content.Find.Execute("135×2=270", $true, $false, $false, $false, $false, $true, 1, $false, "953×4=3812", 2) | Out-Null
$d.Content.Find.Execute("193×9=1737", $true, $false, $false, $false, $false, $true, 1, $false, "967×5=4835", 2) | Out-Null
$d.Content.Find.Execute("756×5=3780", $true, $false, $false, $false, $false, $true, 1, $false, "268×5=1340", 2) | Out-Null
$d.Content.Find.Execute("463×2=926", $true, $false, $false, $false, $false, $true, 1, $false, "743×5=3715", 2) | Out-Null
$d.Content.Find.Execute("391×9=3519", $true, $false, $false, $false, $false, $true, 1, $false, "864×4=3456", 2) | Out-Null
$d.Content.Find.Execute("475×8=3800", $true, $false, $false, $false, $false, $true, 1, $false, "642×8=5136", 2) | Out-Null
$d.Content.Find.Execute("949×9=8541", $true, $false, $false, $false, $false, $true, 1, $false, "899×6=5394", 2) | Out-Null
$d.Content.Find.Execute("890×9=8010", $true, $false, $false, $false, $false, $true, 1, $false, "320×3=960", 2) | Out-Null
$d.Content.Find.Execute("817×5=4085", $true, $false, $false, $false, $false, $true, 1, $false, "653×5=3265", 2) | Out-Null
$d.Content.Find.Execute("250×4=1000", $true, $false, $false, $false, $false, $true, 1, $false, "137×2=274", 2) | Out-Null
$d.Content.Find.Execute("166×8=1328", $true, $false, $false, $false, $false, $true, 1, $false, "633×7=4431", 2) | Out-Null
$d.Content.Find.Execute("814×7=5698", $true, $false, $false, $false, $false, $true, 1, $false, "947×8=7576", 2) | Out-Null
$d.Content.Find.Execute("501×5=2505", $true, $false, $false, $false, $false, $true, 1, $false, "959×9=8631", 2) | Out-Null
$d.Content.Find.Execute("546×9=4914", $true, $false, $false, $false, $false, $true, 1, $false, "767×5=3835", 2) | Out-Null
$d.Content.Find.Execute("689×2=1378", $true, $false, $false, $false, $false, $true, 1, $false, "775×2=1550", 2) | Out-Null
$d.Content.Find.Execute("946×3=2838", $true, $false, $false, $false, $false, $true, 1, $false, "419×7=2933", 2) | Out-Null
$d.Content.Find.Execute("244×6=1464", $true, $false, $false, $false, $false, $true, 1, $false, "372×8=2976", 2) | Out-Null
